$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Cells.Item(2, 4)
$cell.NumberFormat = "@"
$cell.Value = '42.912.77'
$cell.Style = "Normal"
$ws.Cells.Item(2, 5).Value = '  -1.59%  '

$cell = $ws.Cells.Item(3, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.243.32'
$cell.Style = "Normal"
$ws.Cells.Item(3, 5).Value = '  -2.01%  '

$cell = $ws.Cells.Item(4, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.01'
$cell.Style = "Normal"
$ws.Cells.Item(4, 5).Value = '  +0.36%  '

$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = '115.87'
$cell.Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  +0.46%  '

$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = "@"
$cell.Value = '298.11'
$cell.Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  +12.17%  '

$cell = $ws.Cells.Item(7, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.632'
$cell.Style = "Normal"
$ws.Cells.Item(7, 5).Value = '  -1.99%  '

$cell = $ws.Cells.Item(8, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.01'
$cell.Style = "Normal"
$ws.Cells.Item(8, 5).Value = '  +0.18%  '

$cell = $ws.Cells.Item(9, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.620'
$cell.Style = "Normal"
$ws.Cells.Item(9, 5).Value = '  +0.84%  '

$cell = $ws.Cells.Item(10, 4)
$cell.NumberFormat = "@"
$cell.Value = '46.39'
$cell.Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  -2.28%  '

$cell = $ws.Cells.Item(11, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.0934'
$cell.Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  -0.79%  '

$cell = $ws.Cells.Item(12, 4)
$cell.NumberFormat = "@"
$cell.Value = '9.09'
$cell.Style = "Normal"
$ws.Cells.Item(12, 5).Value = '  -0.94%  '

$ws.Cells.Item(13, 5).Value = '  -2.83%  '

$cell = $ws.Cells.Item(14, 4)
$cell.NumberFormat = "@"
$cell.Value = '15.43'
$cell.Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  +0.12%  '

$cell = $ws.Cells.Item(15, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.891'
$cell.Style = "Normal"
$ws.Cells.Item(15, 5).Value = '  +2.13%  '

$cell = $ws.Cells.Item(16, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.585.58'
$cell.Style = "Normal"
$ws.Cells.Item(16, 5).Value = '  -1.88%  '

$cell = $ws.Cells.Item(17, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.253.40'
$cell.Style = "Normal"
$ws.Cells.Item(17, 5).Value = '  -1.66%  '

$cell = $ws.Cells.Item(18, 4)
$cell.NumberFormat = "@"
$cell.Value = '42.797.13'
$cell.Style = "Normal"
$ws.Cells.Item(18, 5).Value = '  -1.95%  '

$ws.Cells.Item(19, 5).Value = '  +11.63%  '

$ws.Cells.Item(20, 5).Value = '  -2.04%  '

$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = "@"
$cell.Value = '74.16'
$cell.Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  +2.27%  '

$cell = $ws.Cells.Item(22, 4)
$cell.NumberFormat = "@"
$cell.Value = '3.49'
$cell.Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  +21.18%  '

$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.36'
$cell.Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  -3.58%  '

$cell = $ws.Cells.Item(24, 4)
$cell.NumberFormat = "@"
$cell.Value = '232.88'
$cell.Style = "Normal"
$ws.Cells.Item(24, 5).Value = '  -1.64%  '

$cell = $ws.Cells.Item(25, 4)
$cell.NumberFormat = "@"
$cell.Value = '9.40'
$cell.Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  -0.85%  '

$cell = $ws.Cells.Item(26, 4)
$cell.NumberFormat = "@"
$cell.Value = '12.24'
$cell.Style = "Normal"
$ws.Cells.Item(26, 5).Value = '  +5.82%  '

$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.999'
$cell.Style = "Normal"
$ws.Cells.Item(27, 5).Value = '  -1.84%  '

$cell = $ws.Cells.Item(28, 4)
$cell.NumberFormat = "@"
$cell.Value = '40.30'
$cell.Style = "Normal"
$ws.Cells.Item(28, 5).Value = '  -4.01%  '

$ws.Cells.Item(29, 5).Value = '  -1.02%  '

$cell = $ws.Cells.Item(30, 4)
$cell.NumberFormat = "@"
$cell.Value = '3.27'
$cell.Style = "Normal"
$ws.Cells.Item(30, 5).Value = '  -3.34%  '

$cell = $ws.Cells.Item(31, 4)
$cell.NumberFormat = "@"
$cell.Value = '175.80'
$cell.Style = "Normal"
$ws.Cells.Item(31, 5).Value = '  +1.10%  '

$cell = $ws.Cells.Item(32, 4)
$cell.NumberFormat = "@"
$cell.Value = '21.29'
$cell.Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  -2.16%  '

$cell = $ws.Cells.Item(33, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.0910'
$cell.Style = "Normal"
$ws.Cells.Item(33, 5).Value = '  -0.06%  '

$cell = $ws.Cells.Item(34, 4)
$cell.NumberFormat = "@"
$cell.Value = '4.60'
$cell.Style = "Normal"
$ws.Cells.Item(34, 5).Value = '  +16.34%  '

$cell = $ws.Cells.Item(35, 4)
$cell.NumberFormat = "@"
$cell.Value = '5.65'
$cell.Style = "Normal"

$cell = $ws.Cells.Item(36, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.128'
$cell.Style = "Normal"
$ws.Cells.Item(36, 5).Value = '  -1.71%  '

$cell = $ws.Cells.Item(37, 4)
$cell.NumberFormat = "@"
$cell.Value = '4.75'
$cell.Style = "Normal"
$ws.Cells.Item(37, 5).Value = '  +1.23%  '

$ws.Cells.Item(38, 5).Value = '  -1.90%  '

$cell = $ws.Cells.Item(39, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.106'
$cell.Style = "Normal"
$ws.Cells.Item(39, 5).Value = '  +0.15%  '

$cell = $ws.Cells.Item(40, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.61'
$cell.Style = "Normal"
$ws.Cells.Item(40, 5).Value = '  +2.10%  '

$cell = $ws.Cells.Item(41, 4)
$cell.NumberFormat = "@"
$cell.Value = '72.55'
$cell.Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  -2.41%  '

$cell = $ws.Cells.Item(42, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.238'
$cell.Style = "Normal"
$ws.Cells.Item(42, 5).Value = '  +1.04%  '

$cell = $ws.Cells.Item(43, 4)
$cell.NumberFormat = "@"
$cell.Value = '13.52'
$cell.Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  -6.13%  '

$ws.Cells.Item(44, 5).Value = '  -0.09%  '

$ws.Cells.Item(45, 5).Value = '  -2.01%  '

$cell = $ws.Cells.Item(46, 4)
$cell.NumberFormat = "@"
$cell.Value = '5.60'
$cell.Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  -6.53%  '

$ws.Cells.Item(47, 5).Value = '  +4.41%  '

$cell = $ws.Cells.Item(48, 4)
$cell.NumberFormat = "@"
$cell.Value = '107.63'
$cell.Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  +6.83%  '

$cell = $ws.Cells.Item(49, 4)
$cell.NumberFormat = "@"
$cell.Value = '8.61'
$cell.Style = "Normal"
$ws.Cells.Item(49, 5).Value = '  +0.16%  '

$cell = $ws.Cells.Item(50, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.0990'
$cell.Style = "Normal"
$ws.Cells.Item(50, 5).Value = '  -0.98%  '

$cell = $ws.Cells.Item(51, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.471'
$cell.Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  +5.49%  '
